$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap "Santa Lucia" (row 189) and "Belice" (row 190) rows ---
# Before: row189 = Santa Lucia (18,0,15,3,0,0,0) ; row190 = Belice (18,0,13,3,1,0,2)
# After:  row189 = Belice      (18,0,13,3,1,0,2) ; row190 = Santa Lucia (18,0,15,3,0,0,0)
$ws.Cells.Item(189, 1).Value = "Belice"
$ws.Cells.Item(189, 2).Value = 18
$ws.Cells.Item(189, 3).Value = 0
$ws.Cells.Item(189, 4).Value = 13
$ws.Cells.Item(189, 5).Value = 3
$ws.Cells.Item(189, 6).Value = 1
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 2

$ws.Cells.Item(190, 1).Value = "Santa Lucia"
$ws.Cells.Item(190, 2).Value = 18
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 15
$ws.Cells.Item(190, 5).Value = 3
$ws.Cells.Item(190, 6).Value = 0
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 0

# --- Swap "San Cristobal y Nieves" (row 198) and "Burundi" (row 199) rows ---
# Before: row198 = San Cristobal y Nieves (15,0,8,7,0,0,0) ; row199 = Burundi (15,0,7,7,0,0,1)
# After:  row198 = Burundi (15,0,7,7,0,0,1) ; row199 = San Cristobal y Nieves (15,0,8,7,0,0,0)
$ws.Cells.Item(198, 1).Value = "Burundi"
$ws.Cells.Item(198, 2).Value = 15
$ws.Cells.Item(198, 3).Value = 0
$ws.Cells.Item(198, 4).Value = 7
$ws.Cells.Item(198, 5).Value = 7
$ws.Cells.Item(198, 6).Value = 0
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 1

$ws.Cells.Item(199, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(199, 2).Value = 15
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 8
$ws.Cells.Item(199, 5).Value = 7
$ws.Cells.Item(199, 6).Value = 0
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 0

# --- Update country totals (Estados Unidos, row 4) ---
$ws.Range("B4").Value = 1211838
$ws.Range("C4").Value = 23716
$ws.Range("D4").Value = 187036
$ws.Range("E4").Value = 955093
$ws.Range("G4").Value = 1112
$ws.Range("H4").Value = 69709

# --- Update country totals (Argentina, row 57) ---
$ws.Range("B57").Value = 4887
$ws.Range("C57").Value = 104
$ws.Range("E57").Value = 3185
$ws.Range("G57").Value = 14
$ws.Range("H57").Value = 260

# --- Update country totals (Uruguay, row 110) ---
$ws.Range("B110").Value = 657
$ws.Range("C110").Value = 2
$ws.Range("D110").Value = 447
$ws.Range("E110").Value = 193

# --- Update the "last updated" timestamp text (row 1 / A1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 02:03"
